# Commit: add `$.num.` prefix for avoiding the string value in json of
# decimal value from database; upgrade mysql2 and typeorm.
#
# The template sheet had two placeholder tokens that fed numeric totals
# into the exported JSON as plain strings:
#   - D6 "$subTotal.cashInTxn" -> now namespaced "$.num.subTotal.cashInTxn"
#     so the renderer emits a real (non-string) numeric value.
#   - D7 "$cashInTxnTotal" placeholder is replaced outright by a live
#     worksheet formula that sums the "Txn" column for the visible rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "`$.num.subTotal.cashInTxn"
$ws.Range("D7").Formula = "=SUM(OFFSET(D`$5,0,0,ROW()-5,1))"

[void]$ws.Range("E19").Select()
